# SBOM.xlsx update: add "Cost model" and "Provenance" columns to the
# Java and Python dependency tables, rename the "Licence type"/"Licence"
# headers to "Licence"/"Official licence link", and add a "Cost model"
# column to the "External dependencies" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Constants (Excel enumeration values used below)
# ---------------------------------------------------------------------
$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138
$xlCenter = -4108
$xlLeft = -4131
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlPasteFormats = -4122
$xlPasteAll = -4104

# ---------------------------------------------------------------------
# 1. Widen the table: extend headers / banner rows from column E/G out
#    to column G, and introduce the two new columns H ("Cost model")
#    and I ("Provenance").
# ---------------------------------------------------------------------

# Re-merge the section banner rows so they span the new width.
$ws.Range("B6:E6").UnMerge()
$ws.Range("B6:G6").Merge()

$ws.Range("B15:E15").UnMerge()
$ws.Range("B15:G15").Merge()

# Copy the banner formatting (grey fill / borders) across the newly
# added columns F, G (already part of the merge) and H, I (kept as
# separate un-merged cells, matching the diff).
$ws.Range("D6").Copy() | Out-Null
$ws.Range("F6:G6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D15").Copy() | Out-Null
$ws.Range("F15:G15").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C6").Copy() | Out-Null
$ws.Range("H6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("G6").Copy() | Out-Null
$ws.Range("I6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G15").Copy() | Out-Null
$ws.Range("I15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I6").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I6").Borders.Item($xlEdgeRight).Weight = $xlMedium
$ws.Range("I15").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I15").Borders.Item($xlEdgeRight).Weight = $xlMedium

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Java dependencies table (rows 7-13): relabel headers, add the new
#    columns, and widen rows so the new wrapped text fits.
# ---------------------------------------------------------------------

# Header row 7 - copy the "Licence"(G7) header style onto the new F7/G7
# positions, then overwrite text; add H7/I7 headers.
$ws.Range("G7").Copy() | Out-Null
$ws.Range("F7").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F7").Value = "Licence"
$ws.Range("G7").Value = "Official licence link"
$ws.Range("H7").Value = "Cost model"
$ws.Range("I7").Value = "Provenance"

$ws.Range("G6").Copy() | Out-Null
$ws.Range("H7:I7").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("H7").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("H7").Borders.Item($xlEdgeTop).Weight = $xlMedium
$ws.Range("H7").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("H7").Borders.Item($xlEdgeBottom).Weight = $xlThin
$ws.Range("H7").Borders.Item($xlEdgeRight).LineStyle = 0
$ws.Range("H7").Interior.Pattern = 0
$ws.Range("H7").Font.Bold = $true
$ws.Range("I7").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("I7").Borders.Item($xlEdgeTop).Weight = $xlMedium
$ws.Range("I7").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("I7").Borders.Item($xlEdgeBottom).Weight = $xlThin
$ws.Range("I7").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I7").Borders.Item($xlEdgeRight).Weight = $xlMedium
$ws.Range("I7").Interior.Pattern = 0
$ws.Range("I7").Font.Bold = $true

# Data rows 8-13: fill H (Cost model) + I (Provenance); base the style
# on the neighbouring G cell (plain bordered cell) with wrap text, and
# raise the row height to fit two-line content.
$javaRows = 8,9,10,11,12,13
foreach ($r in $javaRows) {
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r`:I$r").PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = 0
    $ws.Range("H$r").WrapText = $true
    $ws.Range("I$r").WrapText = $true
    $ws.Rows.Item($r).RowHeight = 40
}
# Re-apply the heavier bottom border on the last (thick-bordered) row.
$ws.Range("B13:I13").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("B13:I13").Borders.Item($xlEdgeBottom).Weight = $xlMedium
$ws.Range("I13").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I13").Borders.Item($xlEdgeRight).Weight = $xlMedium

$ws.Range("H8").Value = "Free to use and distribute, as per licence"
$ws.Range("I8").Value = "FasterXML LLC https://github.com/FasterXML/"
$ws.Range("H9").Value = "Free to use and distribute, as per licence"
$ws.Range("I9").Value = "FasterXML LLC https://github.com/FasterXML/"
$ws.Range("H10").Value = "Free to use, but not distribute, as per licence - this is a testing framework and as such will not be distributed."
$ws.Range("I10").Value = "junit.org"
$ws.Range("H11").Value = "Free to use and distribute, as per licence"
$ws.Range("I11").Value = "The Apache Software Foundation https://www.apache.org/licenses/"
$ws.Range("H12").Value = "Free to use and distribute, as per licence"
$ws.Range("I12").Value = "The Apache Software Foundation https://www.apache.org/licenses/"
$ws.Range("H13").Value = "Free to use and distribute, as per licence"
$ws.Range("I13").Value = "The Apache Software Foundation https://www.apache.org/licenses/"

# Small spacer row 14 grows a little (keeps H14/I14 blank but styled).
$ws.Range("D14").Copy() | Out-Null
$ws.Range("H14:I14").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Rows.Item(14).RowHeight = 27

# ---------------------------------------------------------------------
# 3. Python dependencies table (rows 16-22)
# ---------------------------------------------------------------------
$ws.Range("G16").Copy() | Out-Null
$ws.Range("H16:I16").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("H16").Borders.Item($xlEdgeRight).LineStyle = 0
$ws.Range("H16").Interior.Pattern = 0
$ws.Range("H16").Font.Bold = $true
$ws.Range("I16").Interior.Pattern = 0
$ws.Range("I16").Font.Bold = $true
$ws.Range("I16").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I16").Borders.Item($xlEdgeRight).Weight = $xlMedium

$ws.Range("F16").Value = "Licence"
$ws.Range("G16").Value = "Official licence link"
$ws.Range("H16").Value = "Cost model"
$ws.Range("I16").Value = "Provenance"

# Merge C:D for each dependency row (17-22), matching the new layout
# where the "Dependency" label spans both columns.
$ws.Range("C16:D16").Merge()

$pyRows = 17,18,19,20,21,22
foreach ($r in $pyRows) {
    $ws.Range("C$r`:D$r").Merge()
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r`:I$r").PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = 0
}
$ws.Range("B22:I22").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("B22:I22").Borders.Item($xlEdgeBottom).Weight = $xlMedium
$ws.Range("I22").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("I22").Borders.Item($xlEdgeRight).Weight = $xlMedium

$ws.Range("H17").Value = "Free to use and distribute, as per licence"
$ws.Range("I17").Value = "https://pypi.org/project/matplotlib/"
$ws.Range("H18").Value = "Free to use and distribute, as per licence"
$ws.Range("I18").Value = "https://pypi.org/project/numpy/"
$ws.Range("H19").Value = "Free to use and distribute, as per licence"
$ws.Range("I19").Value = "https://pypi.org/project/pandas/"
$ws.Range("H20").Value = "Free to use and distribute, as per licence"
$ws.Range("I20").Value = "https://pypi.org/project/seaborn/"
$ws.Range("H21").Value = "Free to use and distribute, as per licence"
$ws.Range("I21").Value = "https://pypi.org/project/requests/"
$ws.Range("H22").Value = "Free to use and distribute, as per licence"
$ws.Range("I22").Value = "https://pypi.org/project/pytest/"

# ---------------------------------------------------------------------
# 4. External dependencies table (rows 24-26): add a "Cost model" column.
# ---------------------------------------------------------------------
$ws.Range("B24:C24").UnMerge()
$ws.Range("B24:D24").Merge()

$ws.Range("C25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D25").Value = "Cost model"

$ws.Range("C26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("D26").Value = "Free up to 10k grants per day"

# ---------------------------------------------------------------------
# 5. Column widths (new + adjusted columns)
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 10.6328125   # E
$ws.Columns.Item(7).ColumnWidth = 35           # G (narrower now)
$ws.Columns.Item(8).ColumnWidth = 44.36328125  # H - new
$ws.Columns.Item(9).ColumnWidth = 32.90625     # I - new

# ---------------------------------------------------------------------
# 6. Selection + dimension bookkeeping (cosmetic, mirrors the diff)
# ---------------------------------------------------------------------
$ws.Range("F26").Select()

$excel.CutCopyMode = 0
